$d = $word.ActiveDocument

# 1) Extend the final existing paragraph ("families.") with a line break
#    followed by new text, exactly as in the authored diff.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter("`vFor the initial implementation and experiments required in this lab, the main focus will be on")

# 2) Append the remaining new paragraphs (each on its own line) after that,
#    using a single InsertAfter call with embedded paragraph marks ("`r")
#    for speed/robustness, then verify the resulting paragraph count.
$newParagraphs = @(
    "CICIDS2017 as the primary dataset, because it offers:",
    "- Rich flow-level features suitable for feature selection and pruning. - Multiple attack types",
    "and realistic background traffic. - Widely used benchmarks that allow comparison with existing",
    "methods.",
    "1.2.2. Preprocessing and Data Splits",
    "Network traffic in these datasets is already transformed into flow records. Each flow",
    "corresponds to an ordered tuple of packets between a source and destination, and is described",
    "by multiple numerical and categorical features (e.g., duration, packet counts, bytes, flags,",
    "statistics over inter-arrival times).",
    "The preprocessing steps are:",
    "1. Filtering and cleaning. - Remove duplicated rows and flows with missing or inconsistent",
    "labels. - Optionally remove extremely rare classes to reduce class fragmentation if needed.",
    "2. Label encoding. - The label space Y includes at least two categories: - Normal traffic -",
    "Attack traffic (which may be further subdivided into attack families such as DoS, DDoS,",
    "PortScan, Botnet, etc.). - For binary classification experiments, all attack types are merged",
    "into a single `"Attack`" class. For multi-class experiments, attack families are kept separate.",
    "3. Train/validation/test splitting. - Use time-aware splits where possible, to better simulate",
    "deployment: earlier flows for training and later flows for validation/test. - Alternatively,",
    "use stratified splits to preserve the class distribution, e.g., 60% train, 20% validation, 20%",
    "test.",
    "4. Normalization / scaling. - Continuous features are standardized or normalized (e.g., min–max",
    "scaling) as required by the model family (especially for logistic regression and neural",
    "models). - Tree-based models (Random Forest) may be trained without explicit scaling.",
    "5. Class imbalance handling. - To address the frequent imbalance between normal and attack",
    "traffic, techniques like class-weighting, undersampling of the majority class, or oversampling",
    "of minority attack classes may be used.",
    "1.3. Feature Space and Mathematical Representation",
    "Each network flow is represented by a feature vector in a d-dimensional space.",
    "Let x_i in R^d denote the feature vector for the i-th flow. Let y_i in Y denote the",
    "corresponding class label, where: - For binary classification, Y = {0, 1} with 0 = `"Normal`", 1",
    "= `"Attack`". - For multi-class classification, Y = {0, 1, …, K-1} with K distinct traffic types",
    "(one normal + multiple attack families).",
    "The full dataset is therefore:",
    "D = { (x_i, y_i) } for i = 1..N,",
    "where N is the total number of flows.",
    "The research distinguishes between:",
    "- A full feature set F_full of size d_full, e.g., all flow features provided by CICIDS2017",
    "(typically 70–80 features). - One or more compact feature sets F_compact^(m) of size d_m, where",
    "d_m << d_full (e.g., 20, 30, or 40 features), selected using feature importance methods and",
    "domain knowledge.",
    "Formally, for a given compact set F_compact, the feature mapping reduces each original flow",
    "vector x_i to a lower-dimensional representation:",
    "z_i = phi(x_i) in R^{d_compact},",
    "where phi selects only the coordinates corresponding to the chosen subset of features.",
    "The central mathematical problem is to learn a classification function:",
    "f_theta: R^d -> Y,",
    "parametrized by theta, such that f_theta(x_i) predicts the correct label y_i for unseen flows,",
    "while the computation of f_theta remains efficient on constrained hardware.",
    "1.4. Model Families and Baseline Approaches",
    "To test the hypotheses and to justify the advantages of the proposed lightweight approach, the",
    "experiments will compare several model families:",
    "1. Logistic Regression (LR). - A linear model used as a simple, interpretable baseline. - Works",
    "with a compact feature set and allows fast inference and straightforward feature importance",
    "analysis (via coefficients).",
    "2. Random Forest (RF) / Extremely Randomized Trees (ExtraTrees). - Ensemble tree methods that",
    "handle non-linear relationships and interactions between features. - Provide feature importance",
    "scores and typically achieve strong performance on flow-based intrusion detection tasks.",
    "3. Gradient Boosting (e.g., XGBoost / LightGBM) – if time permits. - More advanced tree",
    "ensembles that may improve accuracy at the cost of increased complexity. - Useful for exploring",
    "the upper bound of performance on the chosen feature sets.",
    "4. Reference Deep / Autoencoder-based methods (baseline from literature). - Although the main",
    "focus is on lightweight models, key results from deep or autoencoder-based intrusion detection",
    "(e.g., Kitsune, N-BaIoT-style models) will be considered as reference baselines from the",
    "literature. - Wherever runtime measurements are available, they will be used to position the",
    "proposed lightweight models in the accuracy–latency–memory space.",
    "For each model f_theta, the training procedure seeks parameters theta* that minimize an",
    "empirical loss on the training set, for example the cross-entropy loss in the multi-class case:",
    "L(theta) = (1 / N_train) * sum over i in Train of l(f_theta(x_i), y_i),",
    "where l is the cross-entropy between the predicted class probabilities and the ground truth",
    "labels.",
    "1.5. Experimental Scenarios",
    "To systematically evaluate the proposed approach, the experiments are structured into several",
    "scenarios, each addressing part of the research questions.",
    "Scenario S1: Full vs. Compact Feature Sets",
    "Goal: Quantify the loss (or possible improvement) in detection performance when moving from the",
    "full feature set to compact subsets. Setup: - Train the same model family (e.g., Random Forest,",
    "logistic regression) on: - Full feature vector x_i in R^{d_full}, - Compact feature vector z_i",
    "in R^{d_compact} with different sizes (e.g., 20, 30, 40 features). - Evaluate on the same",
    "validation and test splits. Expected Outcome: - Identify minimal feature sets that achieve",
    "performance within 1–2% of the full model.",
    "Scenario S2: Comparison of Model Families",
    "Goal: Compare different model families in terms of accuracy and deployment-oriented metrics",
    "using the same compact feature set. Setup: - Fix a compact feature set F_compact (e.g., 30",
    "features). - Train and evaluate LR, Random Forest, ExtraTrees, and (optionally) a gradient",
    "boosting model. - Compare their F1, ROC-AUC, confusion matrices, and runtime measurements.",
    "Expected Outcome: - Determine which model family offers the best trade-off between performance",
    "and efficiency.",
    "Scenario S3: Binary vs. Multi-Class Detection",
    "Goal: Understand the difference in performance and complexity between binary detection (normal",
    "vs. attack) and multi-class classification (distinguishing attack families). Setup: - Train",
    "models for binary classification and for multi-class classification using the same data and"
)

$joined = [string]::Join("`r", $newParagraphs)
$tailRange = $d.Paragraphs.Last.Range
$tailRange.Collapse(0)
$tailRange.InsertAfter("`r" + $joined)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
